# Adds a new "2020" column (N) to the first worksheet, mirroring the
# formatting of the existing 2019 column (M), and moves the selection to N9.
#
# (The source workbook's absPath bookkeeping hint under
# mc:AlternateContent/x15ac:absPath - a purely cosmetic, Excel-managed record
# of the folder the file was last saved from on the author's machine - has no
# corresponding property on the Excel object model, so it cannot be touched
# from automation code; it is left as-is.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header year in N4, formatted like the existing D4:M4 year cells.
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2020

# New data point in N5, formatted like the existing D5:M5 value cells.
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 2.1

$excel.CutCopyMode = $false

# Move the active selection to N9, matching the saved view state.
$ws.Range("N9").Select()
